# "added results for sex" — append 4 new result rows (18-21) to the
# Test_results table on Ark1 for the new "Sex" classifier, matching the
# existing table's layout: Program | Essay | Classifier | Percentages | n-gram

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: okcupid.py / Essay 4 / Sex / 0.56 / Unigram
$ws.Cells.Item(18, 1).Value = "okcupid.py"
$ws.Cells.Item(18, 2).Value = "Essay 4"
$ws.Cells.Item(18, 3).Value = "Sex"
$ws.Cells.Item(18, 4).Value = 0.56
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 5).Value = "Unigram"

# Row 19: okcupid.py / Essay 4 / Sex / 0.58 / Bigram
$ws.Cells.Item(19, 1).Value = "okcupid.py"
$ws.Cells.Item(19, 2).Value = "Essay 4"
$ws.Cells.Item(19, 3).Value = "Sex"
$ws.Cells.Item(19, 4).Value = 0.58
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 5).Value = "Bigram"

# Row 20: okcupid.py / Essay 7 / Sex / 0.58 / Unigram
$ws.Cells.Item(20, 1).Value = "okcupid.py"
$ws.Cells.Item(20, 2).Value = "Essay 7"
$ws.Cells.Item(20, 3).Value = "Sex"
$ws.Cells.Item(20, 4).Value = 0.58
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 5).Value = "Unigram"

# Row 21: okcupid.py / Essay 7 / Sex / 0.59 / Bigram
$ws.Cells.Item(21, 1).Value = "okcupid.py"
$ws.Cells.Item(21, 2).Value = "Essay 7"
$ws.Cells.Item(21, 3).Value = "Sex"
$ws.Cells.Item(21, 4).Value = 0.59
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 5).Value = "Bigram"

# Final selection left on C19, matching the saved workbook state
$ws.Range("C19").Select()
